$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 179, shifting the existing rows 179-286 down to 180-287.
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row 179 with the new weekly record.
$ws.Cells.Item(179, 1).Value = 8
$ws.Cells.Item(179, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(179, 3).Value = "Coquimbo"
$ws.Cells.Item(179, 4).Value = (Get-Date -Year 2022 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(179, 5).Value = 4
$ws.Cells.Item(179, 6).Value = 100114013
$ws.Cells.Item(179, 7).Value = "Zanahoria"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 600
$ws.Cells.Item(179, 11).Value = 5500
$ws.Cells.Item(179, 12).Value = 6000
$ws.Cells.Item(179, 13).Value = 5750
$ws.Cells.Item(179, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(179, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(179, 16).Value = 288
$ws.Cells.Item(179, 17).Value = 20
$ws.Cells.Item(179, 18).Value = "Hortaliza"
